$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.448.36"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.16%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.669.46"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.78%  "

# Row 4
$ws.Range("E4").Value = "  +0.21%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.66"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.94%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5260"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.09%  "

# Row 7
$ws.Range("E7").Value = "  +0.17%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2665"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.57%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06367"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.72%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.63"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.11%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07802"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.67%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.685.34"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.83%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.465"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.25%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5531"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.02%  "

# Row 15
$ws.Range("D15").Value = "0.0₅8269"
$ws.Range("E15").Value = "  +1.01%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.40"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.58%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.460.89"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.19%  "

# Row 18
$ws.Range("E18").Value = "  +0.08%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.734"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.08%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.20"
$ws.Range("D20").ClearFormats()

# Row 21
$ws.Range("E21").Value = "  +2.06%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.264"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.28%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.005"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.32%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1261"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.94%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.62"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.32%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.403"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.72%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.21"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.28%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.418"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.18%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06162"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.33%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.286"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.44%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.613"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +6.61%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.389"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.11%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.677"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.61%  "

# Row 34
$ws.Range("E34").Value = "  +1.81%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6061"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +8.43%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.418"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.97%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.769"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.26%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01610"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.83%  "

# Row 39
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.025"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.33%  "

# Row 40
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.090.09"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +7.19%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8575"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.90%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.002"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.05%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.64"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.18%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.813.40"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.40%  "

# Row 45
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₈110"
$ws.Range("E45").Value = "  -0.46%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.99"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.73%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.150"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.54%  "

# Row 48
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.003"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.11%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05205"
$ws.Range("D49").ClearFormats()

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.480"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +8.03%  "

# Row 51
$ws.Range("E51").Value = "  +0.58%  "
